# GL update - Item build report
# - changed the template of exported file
# - added as_of and period in custom filtering
#
# The report template shrank from 8 rows to 6: the old thick-bottom-border
# separator row (row 4) and the two trailing rows (7-8) are gone, the
# "Out Qty" header label in I5 is cleared (its bold+center formatting is
# kept for later reuse), and a left-aligned placeholder is introduced at G3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the two trailing rows (formerly rows 7 & 8) - sheet shrinks to A1:I6
$ws.Rows("7:8").Delete() | Out-Null

# Remove the medium bottom border (and its thick-bottom row styling) that
# used to separate the header block from the rest of the sheet.
$ws.Range("A4:I4").Borders.LineStyle = -4142   # xlLineStyleNone
$ws.Rows("4:4").AutoFit() | Out-Null           # drop the now-stale 15.75pt/thickBot row height

# New left-aligned formatting placeholder cell for the custom filtering
# (as_of / period) the exported report now carries.
$ws.Range("G3").HorizontalAlignment = -4131    # xlLeft

# The "Out Qty" label text is no longer written into the template here -
# it keeps its bold + centered style, just without the literal text.
$ws.Range("I5").ClearContents() | Out-Null

# Park the active selection past the new used range, like the source file.
$ws.Range("J12").Select() | Out-Null
